# backwardElimination.xlsx - refresh the statsmodels OLS summary timestamps.
# Each of the 16 sheets has its "OLS Regression Results" text block in cell B2;
# update the "Date:" and "Time:" lines from the stale run to the latest run.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("B2")
    $text = $cell.Text
    if ($text -and $text.Contains("Date:")) {
        $newText = $text -replace "Wed, 01 Jan 2020", "Thu, 02 Jan 2020"
        $newText = $newText -replace "23:18:52", "20:48:45"
        $newText = $newText -replace "23:18:53", "20:48:45"
        $cell.Value = $newText
    }
}
